$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "Sender*" / "Destination Facility*" header labels (columns I and J)
$ws.Range("I1").Value = "Destination Facility*"
$ws.Range("J1").Value = "Sender*"

# Update sample data row
$ws.Range("A2").Value = 3
$ws.Range("B2").Value = "'03/19/2022"
$ws.Range("H2").Value = "test11"
$ws.Range("I2").Value = "test12"
$ws.Range("J2").Value = "sender11"
$ws.Range("K2").Value = "rec12"

# Adjust column widths to match new (swapped) best-fit content
$ws.Columns.Item(9).ColumnWidth = 18.33
$ws.Columns.Item(10).ColumnWidth = 10.5

# Update the selection to column K (entire column selected)
$ws.Range("K1:K1048576").Select()

$wb.Save()
